# The deck ships two embedded themes:
#   theme1.xml ("Integral")     -> used by the (only) Slide Master
#   theme2.xml ("Office Theme") -> used by the Notes Master
#
# The author switched the presentation's design from "Integral" to the
# built-in "Office Theme" palette (Design tab -> Themes -> Office Theme).
# That repaints every slide via the Slide Master's color scheme with the
# standard Office theme colors:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

function HexToBgrLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's ColorFormat.RGB is a VBA-style RGB() long: R + G*256 + B*65536
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme color scheme, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToBgrLong $officeThemeHex[$i - 1]
}
